$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-25 Saturday" "2025-10-26 Sunday"
Replace-Text "614÷6=" "153÷6="
Replace-Text "596÷9=" "482÷9="
Replace-Text "468÷2=" "536÷7="
Replace-Text "255÷3=" "659÷6="
Replace-Text "127÷2=" "913÷4="
Replace-Text "591÷4=" "842÷2="
Replace-Text "789÷4=" "660÷6="
Replace-Text "950÷4=" "836÷4="
Replace-Text "387÷9=" "229÷7="
Replace-Text "511÷7=" "122÷9="
Replace-Text "674÷4=" "304÷3="
Replace-Text "180÷8=" "631÷2="
Replace-Text "897÷7=" "479÷2="
Replace-Text "225÷8=" "428÷2="
Replace-Text "793÷8=" "976÷8="
Replace-Text "281÷8=" "748÷2="
Replace-Text "857÷3=" "579÷5="
Replace-Text "868÷5=" "838÷4="
Replace-Text "944÷3=" "735÷9="
Replace-Text "949÷6=" "713÷4="
Replace-Text "194÷7=" "393÷7="
Replace-Text "257÷6=" "281÷6="
Replace-Text "958÷6=" "860÷6="
Replace-Text "259÷3=" "407÷5="
Replace-Text "979÷7=" "429÷7="
